# Updated cryptos list on Fri Feb 17 22:49:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as plain text (e.g. "1.000", "24.539.95").
# Force the whole price column to Text format first so that Excel does not
# silently reinterpret the assigned strings as numbers/dates and strip
# meaningful trailing zeros or thousands-style dot separators.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "24.539.95"
$ws.Range("E2").Value = "  +0.95%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.691.35"
$ws.Range("E3").Value = "  +1.49%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.64%  "

# Row 5 - BNB
$ws.Range("D5").Value = "312.30"
$ws.Range("E5").Value = "  +0.41%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.63%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.3931"
$ws.Range("E7").Value = "  +0.47%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.4026"
$ws.Range("E8").Value = "  +2.03%  "

# Row 9 - was Polygon, now BinanceUSD (rows 9/10 swapped)
$ws.Range("B9").Value = "BinanceUSD"
$ws.Range("C9").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D9").Value = "1.001"
$ws.Range("E9").Value = "  -0.63%  "

# Row 10 - was BinanceUSD, now Polygon
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.518"
$ws.Range("E10").Value = "  +8.64%  "

# Row 11 - OKB
$ws.Range("D11").Value = "53.41"
$ws.Range("E11").Value = "  +9.65%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "0.08754"
$ws.Range("E12").Value = "  +1.83%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "7.296"
$ws.Range("E13").Value = "  +12.56%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "0.00001315"
$ws.Range("E15").Value = "  +2.47%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "7.539"
$ws.Range("E16").Value = "  +6.10%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.693.31"
$ws.Range("E17").Value = "  +1.17%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "100.18"
$ws.Range("E18").Value = "  -0.83%  "

# Row 19 - TRON
$ws.Range("D19").Value = "0.07071"
$ws.Range("E19").Value = "  +4.35%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "19.43"
$ws.Range("E20").Value = "  +2.59%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.692"
$ws.Range("E21").Value = "  +0.09%  "

# Row 22 - Dai
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.58%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "14.12"
$ws.Range("E23").Value = "  +3.36%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "24.534.38"
$ws.Range("E24").Value = "  +0.97%  "

# Row 25 - LidoDAOToken
$ws.Range("D25").Value = "3.028"
$ws.Range("E25").Value = "  +9.74%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "2.311"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "22.28"
$ws.Range("E27").Value = "  +2.17%  "

# Row 28 - Monero
$ws.Range("D28").Value = "159.28"
$ws.Range("E28").Value = "  +0.43%  "

# Row 29 - HuobiToken
$ws.Range("D29").Value = "5.155"
$ws.Range("E29").Value = "  -1.74%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "133.34"
$ws.Range("E30").Value = "  +1.66%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "7.523"
$ws.Range("E31").Value = "  +31.10%  "

# Row 32 - WrappedliquidstakedEther2.0
$ws.Range("D32").Value = "1.880.28"
$ws.Range("E32").Value = "  +1.29%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "1.089"
$ws.Range("E33").Value = "  -3.82%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.08630"
$ws.Range("E34").Value = "  +1.61%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Range("D35").Value = "7.320"
$ws.Range("E35").Value = "  +20.31%  "

# Row 36 - WEMIXTOKEN
$ws.Range("D36").Value = "1.963"
$ws.Range("E36").Value = "  +9.12%  "

# Row 37 - FraxShare
$ws.Range("D37").Value = "11.00"
$ws.Range("E37").Value = "  +5.52%  "

# Row 38 - Algorand
$ws.Range("D38").Value = "0.2708"
$ws.Range("E38").Value = "  +3.66%  "

# Row 39 - Aptos
$ws.Range("D39").Value = "14.70"
$ws.Range("E39").Value = "  -2.31%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "0.02752"
$ws.Range("E40").Value = "  +9.95%  "

# Row 41 - Stellar
$ws.Range("D41").Value = "0.08972"
$ws.Range("E41").Value = "  +2.21%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "1.472"
$ws.Range("E42").Value = "  +3.09%  "

# Row 43 - TheSandbox
$ws.Range("D43").Value = "0.7626"
$ws.Range("E43").Value = "  +4.16%  "

# Row 44 - Decentraland
$ws.Range("D44").Value = "0.7135"
$ws.Range("E44").Value = "  +2.31%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "15.40"
$ws.Range("E45").Value = "  +3.78%  "

# Row 46 - NEARProtocol
$ws.Range("D46").Value = "2.443"
$ws.Range("E46").Value = "  +2.74%  "

# Row 47 - PancakeSwap
$ws.Range("D47").Value = "4.155"
$ws.Range("E47").Value = "  +1.87%  "

# Row 48 - Frax
$ws.Range("D48").Value = "0.9996"
$ws.Range("E48").Value = "  -0.69%  "

# Row 49 - Quant
$ws.Range("D49").Value = "140.22"
$ws.Range("E49").Value = "  +1.14%  "

# Row 50 - was Flow, now BabyDogeCoin (rows 50/51 swapped)
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000381"
$ws.Range("E50").Value = "  +4.57%  "

# Row 51 - was BabyDogeCoin, now Flow
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "1.291"
$ws.Range("E51").Value = "  +14.53%  "
